$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 29
$ws.Range("B8").Value = "Update index.py"
$ws.Range("C8").Value = "riya-morankar"
$ws.Range("D8").Value = "N/A"
$ws.Range("E8").Value = "edit2 to main"
$ws.Range("F8").Value = "'2025-06-17"
$ws.Range("F8").Style = "Normal"
